# Auto-generated edits applying the Halicarnassus_Profits market-price refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 108.933334
$ws.Cells.Item(2, 9).Value = 109.42857
$ws.Cells.Item(2, 10).Value = 102
$ws.Cells.Item(2, 11).Value = 109.42857
$ws.Cells.Item(2, 12).Value = 102
$ws.Cells.Item(2, 13).Value = 3.571430000000007
$ws.Cells.Item(2, 14).Value = -328
$ws.Cells.Item(38, 8).Value = 3394
$ws.Cells.Item(38, 9).Value = 1716.6
$ws.Cells.Item(38, 11).Value = 5149.799999999999
$ws.Cells.Item(38, 13).Value = -4777.799999999999
$ws.Cells.Item(40, 8).Value = 6253.606
$ws.Cells.Item(40, 9).Value = 4743.467
$ws.Cells.Item(40, 11).Value = 4743.467
$ws.Cells.Item(40, 13).Value = -4568.467
$ws.Cells.Item(58, 8).Value = 1862.5
$ws.Cells.Item(58, 10).Value = 2749.1667
$ws.Cells.Item(58, 12).Value = 8247.500100000001
$ws.Cells.Item(58, 14).Value = -8547.500100000001
$ws.Cells.Item(62, 8).Value = 8826.5
$ws.Cells.Item(62, 9).Value = 6710.8335
$ws.Cells.Item(62, 11).Value = 6710.8335
$ws.Cells.Item(62, 13).Value = -6086.8335
$ws.Cells.Item(64, 8).Value = 8666.666999999999
$ws.Cells.Item(64, 10).Value = 8666.666999999999
$ws.Cells.Item(64, 12).Value = 8666.666999999999
$ws.Cells.Item(64, 14).Value = -9162.666999999999
$ws.Cells.Item(65, 8).Value = 8826.5
$ws.Cells.Item(65, 9).Value = 6710.8335
$ws.Cells.Item(65, 11).Value = 33554.1675
$ws.Cells.Item(65, 13).Value = -30434.1675
$ws.Cells.Item(67, 8).Value = 8666.666999999999
$ws.Cells.Item(67, 10).Value = 8666.666999999999
$ws.Cells.Item(67, 12).Value = 8666.666999999999
$ws.Cells.Item(67, 14).Value = -10382.667
$ws.Cells.Item(113, 8).Value = 4392.8423
$ws.Cells.Item(113, 9).Value = 3072.625
$ws.Cells.Item(113, 10).Value = 5353
$ws.Cells.Item(113, 11).Value = 3072.625
$ws.Cells.Item(113, 12).Value = 5353
$ws.Cells.Item(113, 13).Value = 181.375
$ws.Cells.Item(113, 14).Value = -11861
$ws.Cells.Item(116, 8).Value = 6437.25
$ws.Cells.Item(116, 9).Value = 5999.4
$ws.Cells.Item(116, 10).Value = 7167
$ws.Cells.Item(116, 11).Value = 5999.4
$ws.Cells.Item(116, 12).Value = 7167
$ws.Cells.Item(116, 13).Value = -2557.4
$ws.Cells.Item(116, 14).Value = -14051
$ws.Cells.Item(137, 8).Value = 2996.6843
$ws.Cells.Item(137, 9).Value = 1286.6666
$ws.Cells.Item(137, 10).Value = 3317.3125
$ws.Cells.Item(137, 11).Value = 3859.9998
$ws.Cells.Item(137, 12).Value = 9951.9375
$ws.Cells.Item(137, 13).Value = -1309.9998
$ws.Cells.Item(137, 14).Value = -15051.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4004346.2
$ws.Cells.Item(32, 9).Value = 4737.5
$ws.Cells.Item(32, 10).Value = 14289055
$ws.Cells.Item(32, 11).Value = 4737.5
$ws.Cells.Item(32, 12).Value = 14289055
$ws.Cells.Item(32, 13).Value = -4450.5
$ws.Cells.Item(32, 14).Value = -14289629
$ws.Cells.Item(61, 8).Value = 6611.1113
$ws.Cells.Item(61, 9).Value = 5250
$ws.Cells.Item(61, 10).Value = 7700
$ws.Cells.Item(61, 11).Value = 5250
$ws.Cells.Item(61, 12).Value = 7700
$ws.Cells.Item(61, 13).Value = -5038
$ws.Cells.Item(61, 14).Value = -8124
$ws.Cells.Item(74, 8).Value = 2996
$ws.Cells.Item(74, 9).Value = 2996
$ws.Cells.Item(74, 11).Value = 2996
$ws.Cells.Item(74, 13).Value = -2122
$ws.Cells.Item(77, 8).Value = 2996
$ws.Cells.Item(77, 9).Value = 2996
$ws.Cells.Item(77, 11).Value = 14980
$ws.Cells.Item(77, 13).Value = -10612
$ws.Cells.Item(132, 8).Value = 3725
$ws.Cells.Item(132, 9).Value = 3725
$ws.Cells.Item(132, 11).Value = 11175
$ws.Cells.Item(132, 13).Value = -8645
$ws.Cells.Item(136, 8).Value = 6611.1113
$ws.Cells.Item(136, 9).Value = 5250
$ws.Cells.Item(136, 10).Value = 7700
$ws.Cells.Item(136, 11).Value = 15750
$ws.Cells.Item(136, 12).Value = 23100
$ws.Cells.Item(136, 13).Value = -13200
$ws.Cells.Item(136, 14).Value = -28200

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1501.2609
$ws.Cells.Item(58, 9).Value = 1024.5
$ws.Cells.Item(58, 11).Value = 1024.5
$ws.Cells.Item(58, 13).Value = -821.5
$ws.Cells.Item(102, 8).Value = 241999
$ws.Cells.Item(102, 10).Value = 241999
$ws.Cells.Item(102, 12).Value = 241999
$ws.Cells.Item(102, 14).Value = -246867
$ws.Cells.Item(105, 8).Value = 1473.6666
$ws.Cells.Item(105, 9).Value = 883.9
$ws.Cells.Item(105, 11).Value = 883.9
$ws.Cells.Item(105, 13).Value = 863.1
$ws.Cells.Item(136, 8).Value = 1501.2609
$ws.Cells.Item(136, 9).Value = 1024.5
$ws.Cells.Item(136, 11).Value = 3073.5
$ws.Cells.Item(136, 13).Value = -523.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 118.5625
$ws.Cells.Item(12, 9).Value = 6.8
$ws.Cells.Item(12, 10).Value = 169.36363
$ws.Cells.Item(12, 11).Value = 20.4
$ws.Cells.Item(12, 12).Value = 508.09089
$ws.Cells.Item(12, 13).Value = 152.6
$ws.Cells.Item(12, 14).Value = -854.0908899999999
$ws.Cells.Item(68, 8).Value = 3248.2856
$ws.Cells.Item(68, 9).Value = 1915
$ws.Cells.Item(68, 10).Value = 4248.25
$ws.Cells.Item(68, 11).Value = 5745
$ws.Cells.Item(68, 12).Value = 12744.75
$ws.Cells.Item(68, 13).Value = -4934
$ws.Cells.Item(68, 14).Value = -14366.75
$ws.Cells.Item(71, 8).Value = 3248.2856
$ws.Cells.Item(71, 9).Value = 1915
$ws.Cells.Item(71, 10).Value = 4248.25
$ws.Cells.Item(71, 11).Value = 17235
$ws.Cells.Item(71, 12).Value = 38234.25
$ws.Cells.Item(71, 13).Value = -13179
$ws.Cells.Item(71, 14).Value = -46346.25
$ws.Cells.Item(81, 8).Value = 1071.3334
$ws.Cells.Item(81, 9).Value = 200
$ws.Cells.Item(81, 10).Value = 1507
$ws.Cells.Item(81, 11).Value = 600
$ws.Cells.Item(81, 12).Value = 4521
$ws.Cells.Item(81, 13).Value = 523
$ws.Cells.Item(81, 14).Value = -6767
$ws.Cells.Item(84, 8).Value = 1071.3334
$ws.Cells.Item(84, 9).Value = 200
$ws.Cells.Item(84, 10).Value = 1507
$ws.Cells.Item(84, 11).Value = 1800
$ws.Cells.Item(84, 12).Value = 13563
$ws.Cells.Item(84, 13).Value = 3816
$ws.Cells.Item(84, 14).Value = -24795
$ws.Cells.Item(111, 8).Value = 174.5
$ws.Cells.Item(111, 9).Value = 174.5
$ws.Cells.Item(111, 11).Value = 523.5
$ws.Cells.Item(111, 13).Value = 2543.5
$ws.Cells.Item(113, 8).Value = 2211.1538
$ws.Cells.Item(113, 9).Value = 494.42856
$ws.Cells.Item(113, 11).Value = 1483.28568
$ws.Cells.Item(113, 13).Value = 686.71432
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).ClearContents()
$ws.Cells.Item(116, 14).ClearContents()
$ws.Cells.Item(118, 8).Value = 1360
$ws.Cells.Item(118, 9).Value = 720
$ws.Cells.Item(118, 10).Value = 2000
$ws.Cells.Item(118, 11).Value = 2160
$ws.Cells.Item(118, 12).Value = 6000
$ws.Cells.Item(118, 13).Value = -917
$ws.Cells.Item(118, 14).Value = -8486
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value = 42185.375
$ws.Cells.Item(29, 10).Value = 40425.5
$ws.Cells.Item(29, 12).Value = 40425.5
$ws.Cells.Item(29, 14).Value = -41005.5
$ws.Cells.Item(80, 8).Value = 5441.5
$ws.Cells.Item(80, 10).Value = 1000
$ws.Cells.Item(80, 12).Value = 1000
$ws.Cells.Item(80, 14).Value = -2996
$ws.Cells.Item(83, 8).Value = 5441.5
$ws.Cells.Item(83, 10).Value = 1000
$ws.Cells.Item(83, 12).Value = 5000
$ws.Cells.Item(83, 14).Value = -14984
$ws.Cells.Item(132, 8).Value = 4827.143
$ws.Cells.Item(132, 9).Value = 5663.3335
$ws.Cells.Item(132, 11).Value = 16990.0005
$ws.Cells.Item(132, 13).Value = -14460.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(20, 8).Value = 12500
$ws.Cells.Item(20, 9).Value = 5000
$ws.Cells.Item(20, 11).Value = 5000
$ws.Cells.Item(20, 13).Value = -4774
$ws.Cells.Item(42, 8).Value = 29499
$ws.Cells.Item(42, 9).Value = 26248.75
$ws.Cells.Item(42, 10).Value = 35999.5
$ws.Cells.Item(42, 11).Value = 26248.75
$ws.Cells.Item(42, 12).Value = 35999.5
$ws.Cells.Item(42, 13).Value = -25685.75
$ws.Cells.Item(42, 14).Value = -37125.5
$ws.Cells.Item(49, 8).Value = 29499
$ws.Cells.Item(49, 9).Value = 26248.75
$ws.Cells.Item(49, 10).Value = 35999.5
$ws.Cells.Item(49, 11).Value = 26248.75
$ws.Cells.Item(49, 12).Value = 35999.5
$ws.Cells.Item(49, 13).Value = -26101.75
$ws.Cells.Item(49, 14).Value = -36293.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 40964.3
$ws.Cells.Item(4, 9).Value = 40964.3
$ws.Cells.Item(4, 11).Value = 40964.3
$ws.Cells.Item(4, 13).Value = -40851.3

